$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05157966666666666
$ws.Range("M2").Value = 0.00535
$ws.Range("N2").Value = 0.01605
$ws.Range("O2").Value = 0.003591913026022235
$ws.Range("P2").Value = 0.003591913026022235
$ws.Range("Q2").Value = 0.0002759512166666666
$ws.Range("R2").Value = 0.002483560949999999
$ws.Range("S2").Value = 0.003591913026022235
$ws.Range("T2").Value = 0.003591913026022235

$ws.Range("G3").Value = 0.05157966666666666
$ws.Range("O3").Value = 0.9964080869739778
$ws.Range("P3").Value = 0.9964080869739778
$ws.Range("Q3").Value = 0.07654974435766665
$ws.Range("R3").Value = 0.6889476992189999
$ws.Range("S3").Value = 0.9964080869739778
$ws.Range("T3").Value = 0.9964080869739778
